$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14, column P (سعر البيع / selling price) for "سرنجات 3 سم": "2.0000" -> "6.0000"
# The cell stores this as text (shared string) even though its style carries a
# numeric format (0.00), so force it through a Text number-format round-trip to
# stop the engine from auto-coercing the numeric-looking string into a number,
# then restore the original "0.00" number format on the cell.
$ws.Cells.Item(14, 16).NumberFormat = "@"
$ws.Cells.Item(14, 16).Value = "6.0000"
$ws.Cells.Item(14, 16).NumberFormat = "0.00"

# Row 14, column Q (عدد التعاملات / transaction count): "1:0" -> "3:0"
$ws.Cells.Item(14, 17).Value = "3:0"

# Row 17, column P (grand total): 611 -> 615
$ws.Range("P17").Value = 615

# Footer generation timestamp: "Sunday, 13 July, 2025 10:20 AM" -> "...10:44 AM"
$ws.Range("A18").Value = "Sunday, 13 July, 2025 10:44 AM"
